$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly-graded make-up / bonus points (value 5) that were missing
# for several students, plus a couple of other score additions.
$ws.Range("R13").Value = 5

$ws.Range("P18").Value = 5
$ws.Range("R18").Value = 5

$ws.Range("P31").Value = 5
$ws.Range("R31").Value = 5

$ws.Range("P35").Value = 5
$ws.Range("R35").Value = 5

$ws.Range("P40").Value = 5
$ws.Range("R40").Value = 5

$ws.Range("P47").Value = 5
$ws.Range("R47").Value = 5

$ws.Range("H63").Value = 5
$ws.Range("R63").Value = 5

$ws.Range("P68").Value = 5
$ws.Range("R68").Value = 5

$ws.Range("P71").Value = 5

# This student's quiz score was previously a blank placeholder (two spaces);
# now it has an actual numeric score.
$ws.Range("L81").Value = 9

$ws.Range("R88").Value = 5

# Restore the active selection / view position that was recorded when the
# workbook was last saved.
[void]$ws.Range("R31").Select()

$win = $ws.Application.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
